$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 12 CSK vs RR - fill in the raw score values for row 21
$ws.Range("E21").Value = 100
$ws.Range("H21").Value = 80
$ws.Range("K21").Value = 20
$ws.Range("N21").Value = 0
$ws.Range("Q21").Value = 60
$ws.Range("T21").Value = 40
